$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# 1) "Taken verdeling" paragraph: the stray `_GoBack` bookmark that sat here
#    (left over from the previous edit) moves down to the "Sluiting" item below,
#    so drop it from this paragraph while keeping everything else identical.
$rngTaken = $d.Content
$foundTaken = $rngTaken.Find.Execute("Taken verdeling")
if (-not $foundTaken) {
    throw "Could not find 'Taken verdeling' paragraph"
}
$xmlTaken = '<w:p ' + $wns + ' w14:paraId="016D4414" w14:textId="4F58798A" w:rsidR="00DD0255" w:rsidRDefault="005B42F1" w:rsidP="00DD0255">' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Taken verdeling</w:t></w:r>' +
    '</w:p>'
$rngTaken.InsertXML($xmlTaken) | Out-Null

# 2) "9. Sluiting" -> "7. Sluiting": renumbered now that the "Kosten & Baten" and
#    "Risico's" chapters are gone. The edit lands mid-run (splitting "9" / ". Sluiting"),
#    which is why the `_GoBack` bookmark ends up parked right after the "7".
$rngSluiting = $d.Content
$foundSluiting = $rngSluiting.Find.Execute("9. Sluiting")
if (-not $foundSluiting) {
    throw "Could not find '9. Sluiting' paragraph"
}
$rPr = '<w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:color w:val="000000"/><w:lang w:eastAsia="en-US"/></w:rPr>'
$xmlSluiting = '<w:p ' + $wns + ' w14:paraId="00145042" w14:textId="32B61563" w:rsidR="00735478" w:rsidRPr="001542A3" w:rsidRDefault="00735478" w:rsidP="00D312E4">' +
    '<w:pPr>' + $rPr + '</w:pPr>' +
    '<w:r w:rsidRPr="001542A3">' + $rPr + '<w:t>7</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '<w:r>' + $rPr + '<w:t>. Sluiting</w:t></w:r>' +
    '</w:p>'
$rngSluiting.InsertXML($xmlSluiting) | Out-Null
